$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# 1. Insert 4 new columns before the old "Relation_In_In" column (O:R),
#    shifting everything from O onward 4 columns to the right.
$ws.Range("O1:R1").EntireColumn.Insert()

# Give the new columns roughly the same (narrow, default) width as column N,
# matching the formatting Excel applies to freshly inserted columns.
$ws.Range("O1:R1").EntireColumn.ColumnWidth = 9.67

# 2. Populate the new header cells. Order matters for shared-string allocation:
#    ramp_up_Output1, ramp_up_Output2, ramp_down_Output2, ramp_down_Output1
$ws.Range("O1").Value = "ramp_up_Output1"
$ws.Range("P1").Value = "ramp_up_Output2"
$ws.Range("R1").Value = "ramp_down_Output2"
$ws.Range("Q1").Value = "ramp_down_Output1"

# 3. Grow the table ("Table1") to cover the newly inserted columns.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:AC6"))

# Re-apply the names of every column that sits after the insertion point,
# since the simple Resize() call can mis-map them internally.
$lo.HeaderRowRange.Cells(1,19).Value = "Relation_In_In"
$lo.HeaderRowRange.Cells(1,20).Value = "Relation_In_Out"
$lo.HeaderRowRange.Cells(1,21).Value = "Relation_Out_Out"
$lo.HeaderRowRange.Cells(1,22).Value = "Cost_invest"
$lo.HeaderRowRange.Cells(1,23).Value = "unit_on_cost"
$lo.HeaderRowRange.Cells(1,24).Value = "fom_cost"
$lo.HeaderRowRange.Cells(1,25).Value = "vom_cost"
$lo.HeaderRowRange.Cells(1,26).Value = "vom_cost_Input1"
$lo.HeaderRowRange.Cells(1,27).Value = "vom_cost_Input2"
$lo.HeaderRowRange.Cells(1,28).Value = "vom_cost_Output1"
$lo.HeaderRowRange.Cells(1,29).Value = "vom_cost_Output2"

# 4. New ramping data for the "Methanol_Reactor" row (row 6).
$ws.Range("O6").Value = 0.5
$ws.Range("Q6").Value = 0.5

# 5. The minimum unit-on-cost value (previously on row 5, now shifted to W5
#    after the column insert) moves down to row 6 instead. Build it via a
#    helper/TEXT formula and paste-as-values so it stays text (matching the
#    original "0.0000001" shared string) instead of being re-parsed as a
#    number.
$helper = $ws.Range("AZ100")
$helper.Formula = '=TEXT(0.0000001,"0.0000000")'
$helper.Copy()
$ws.Range("W6").PasteSpecial(-4163)  # xlPasteValues
$helper.ClearContents()
$ws.Range("W5").ClearContents()

# 6. Restore the active-cell selection that was recorded for this sheet.
$ws.Range("Q2").Select()
